# Weekly fruit/vegetable price update: a new daily record is inserted
# at the top of the data block (row 161), pushing all subsequent rows
# down by one and extending the used range to row 261.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 161; Excel shifts rows 161:260 down to
# 162:261 and carries the row-161 formatting (e.g. the date style on
# column D) down with them / into the freshly inserted row.
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new weekly record.
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44596
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100112032
$ws.Range("G161").Value = "Zapallo italiano"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 500
$ws.Range("K161").Value = 9500
$ws.Range("L161").Value = 10000
$ws.Range("M161").Value = 9750
$ws.Range("N161").Value = "$/caja 70 unidades"
$ws.Range("O161").Value = "Provincia de Limarí"
$ws.Range("P161").Value = 139
$ws.Range("Q161").Value = 70
$ws.Range("R161").Value = "Hortaliza"
